$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: FAPs / Ccl21b / Cxcr3 / M1
$ws.Range("A2").Value = "FAPs"
$ws.Range("B2").Value = "Ccl21b"
$ws.Range("C2").Value = "Cxcr3"
$ws.Range("D2").Value = "M1"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.126417
$ws.Range("H2").Value = 0.379251
$ws.Range("I2").Value = 0.5842846710605375
$ws.Range("J2").Value = 0.6782743529807329
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.7166990000000001
$ws.Range("N2").Value = 2.150097
$ws.Range("O2").Value = 0.2276207788704612
$ws.Range("P2").Value = 0.2276207788704611
$ws.Range("Q2").Value = 0.09060293748300001
$ws.Range("R2").Value = 0.815426437347
$ws.Range("S2").Value = 0.1329953319088708
$ws.Range("T2").Value = 0.1543893365133325

# Row 3: FAPs / Ccl21b / Cxcr3 / M2
$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Ccl21b"
$ws.Range("C3").Value = "Cxcr3"
$ws.Range("D3").Value = "M2"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.126417
$ws.Range("H3").Value = 0.379251
$ws.Range("I3").Value = 0.5842846710605375
$ws.Range("J3").Value = 0.6782743529807329
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 2.431954666666666
$ws.Range("N3").Value = 7.295864
$ws.Range("O3").Value = 0.7723792211295388
$ws.Range("P3").Value = 0.7723792211295388
$ws.Range("Q3").Value = 0.307440413096
$ws.Range("R3").Value = 2.766963717864
$ws.Range("S3").Value = 0.4512893391516667
$ws.Range("T3").Value = 0.5238850164674004

# Row 4: sCs / Ccl21b / Cxcr3 / M1
$ws.Range("A4").Value = "sCs"
$ws.Range("B4").Value = "Ccl21b"
$ws.Range("C4").Value = "Cxcr3"
$ws.Range("D4").Value = "M1"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.089945
$ws.Range("H4").Value = 0.17989
$ws.Range("I4").Value = 0.4157153289394626
$ws.Range("J4").Value = 0.321725647019267
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.7166990000000001
$ws.Range("N4").Value = 2.150097
$ws.Range("O4").Value = 0.2276207788704612
$ws.Range("P4").Value = 0.2276207788704611
$ws.Range("Q4").Value = 0.064463491555
$ws.Range("R4").Value = 0.38678094933
$ws.Range("S4").Value = 0.09462544696159043
$ws.Range("T4").Value = 0.07323144235712861

# Row 5: sCs / Ccl21b / Cxcr3 / M2
$ws.Range("A5").Value = "sCs"
$ws.Range("B5").Value = "Ccl21b"
$ws.Range("C5").Value = "Cxcr3"
$ws.Range("D5").Value = "M2"
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 0.089945
$ws.Range("H5").Value = 0.17989
$ws.Range("I5").Value = 0.4157153289394626
$ws.Range("J5").Value = 0.321725647019267
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 2.431954666666666
$ws.Range("N5").Value = 7.295864
$ws.Range("O5").Value = 0.7723792211295388
$ws.Range("P5").Value = 0.7723792211295388
$ws.Range("Q5").Value = 0.2187421624933333
$ws.Range("R5").Value = 1.31245297496
$ws.Range("S5").Value = 0.3210898819778721
$ws.Range("T5").Value = 0.2484942046621384
